$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 11.637904816096672
$ws.Range("C2").Value = 11.559030914856949
$ws.Range("D2").Value = 12.618210600673645
$ws.Range("E2").Value = 11.654864853169308

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 11.090405638127915
$ws.Range("C3").Value = 9.6108703042399721
$ws.Range("D3").Value = 11.196423106175704
$ws.Range("E3").Value = 11.341686551876091

# Update selection to match the new authored range
$ws.Range("B1:E3").Select()
